# Update scripts with new TPM: recomputed NATMI ligand/receptor
# expression + specificity statistics for the Serpinf1 -> Plxdc1
# LR-pair sheet after the "ECs" cluster's ligand (Serpinf1) and
# receptor (Plxdc1) expression values were refreshed with new TPM
# figures. The downstream specificity (I,J,O,P) and edge weight /
# specificity (Q,R,S,T) columns are derivatives of the raw
# expression columns, so they change for every row, while the raw
# E:H / M:N columns only change on the rows touching the "ECs"
# cluster (sending rows 2-4, target rows 2/5/8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.432712
$ws.Range("H2").Value = 1.298136
$ws.Range("I2").Value = 0.0008141316526434356
$ws.Range("J2").Value = 0.0008141316526434357
$ws.Range("M2").Value = 3.849718000000001
$ws.Range("N2").Value = 11.549154
$ws.Range("O2").Value = 0.08094976577134179
$ws.Range("P2").Value = 0.08094976577134178
$ws.Range("Q2").Value = 1.665819175216
$ws.Range("R2").Value = 14.992372576944
$ws.Range("S2").Value = 0.0000659037665885215
$ws.Range("T2").Value = 0.0000659037665885215

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.432712
$ws.Range("H3").Value = 1.298136
$ws.Range("I3").Value = 0.0008141316526434356
$ws.Range("J3").Value = 0.0008141316526434357
$ws.Range("O3").Value = 0.1017970208850765
$ws.Range("P3").Value = 0.1017970208850765
$ws.Range("Q3").Value = 2.094822977613333
$ws.Range("R3").Value = 18.85340679852
$ws.Range("S3").Value = 0.00008287617684734566
$ws.Range("T3").Value = 0.00008287617684734566

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.432712
$ws.Range("H4").Value = 1.298136
$ws.Range("I4").Value = 0.0008141316526434356
$ws.Range("J4").Value = 0.0008141316526434357
$ws.Range("O4").Value = 0.8172532133435817
$ws.Range("P4").Value = 0.8172532133435817
$ws.Range("Q4").Value = 16.81778891912
$ws.Range("R4").Value = 151.36010027208
$ws.Range("S4").Value = 0.0006653517092075683
$ws.Range("T4").Value = 0.0006653517092075684

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.9930385075793827
$ws.Range("J5").Value = 0.9930385075793827
$ws.Range("M5").Value = 3.849718000000001
$ws.Range("N5").Value = 11.549154
$ws.Range("O5").Value = 0.08094976577134179
$ws.Range("P5").Value = 0.08094976577134178
$ws.Range("Q5").Value = 2031.885853206242
$ws.Range("R5").Value = 18286.97267885618
$ws.Range("S5").Value = 0.08038623459047385
$ws.Range("T5").Value = 0.08038623459047384

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.9930385075793827
$ws.Range("J6").Value = 0.9930385075793827
$ws.Range("O6").Value = 0.1017970208850765
$ws.Range("P6").Value = 0.1017970208850765
$ws.Range("S6").Value = 0.1010883616957436
$ws.Range("T6").Value = 0.1010883616957436

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.9930385075793827
$ws.Range("J7").Value = 0.9930385075793827
$ws.Range("O7").Value = 0.8172532133435817
$ws.Range("P7").Value = 0.8172532133435817
$ws.Range("S7").Value = 0.8115639112931652
$ws.Range("T7").Value = 0.8115639112931652

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.006147360767973794
$ws.Range("J8").Value = 0.006147360767973795
$ws.Range("M8").Value = 3.849718000000001
$ws.Range("N8").Value = 11.549154
$ws.Range("O8").Value = 0.08094976577134179
$ws.Range("P8").Value = 0.08094976577134178
$ws.Range("Q8").Value = 12.57829911294
$ws.Range("R8").Value = 113.20469201646
$ws.Range("S8").Value = 0.0004976274142794144
$ws.Range("T8").Value = 0.0004976274142794144

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.006147360767973794
$ws.Range("J9").Value = 0.006147360767973795
$ws.Range("O9").Value = 0.1017970208850765
$ws.Range("P9").Value = 0.1017970208850765
$ws.Range("S9").Value = 0.0006257830124855282
$ws.Range("T9").Value = 0.0006257830124855282

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.006147360767973794
$ws.Range("J10").Value = 0.006147360767973795
$ws.Range("O10").Value = 0.8172532133435817
$ws.Range("P10").Value = 0.8172532133435817
$ws.Range("S10").Value = 0.005023950341208851
$ws.Range("T10").Value = 0.005023950341208852
